$d = $word.ActiveDocument
$lq = [char]0x201C
$rq = [char]0x201D

# ---------------------------------------------------------------------------
# Edit 1: "Critical - ... 26a, 531, 1016" -> append a new run ", 1518a"
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(", 531, 1016")
if ($found1) {
  $rng1.Collapse(0)
  $rng1.InsertAfter(", 1518a")
} else {
  Write-Host "WARN: edit 1 target not found"
}

# ---------------------------------------------------------------------------
# Edit 2: "Warnings - 1670" -> remove the run containing "1670"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("1670")
if ($found2) {
  $rng2.Text = ""
} else {
  Write-Host "WARN: edit 2 target not found"
}

# ---------------------------------------------------------------------------
# Edit 3: "353: An error "CRC of the data file is not equal to CRC subfield
# of catalog record"." -> 4 separate runs reading
#   "1016: " / ""CRC " / "values do not match" / "".""
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$target3 = "353: An error " + $lq + "CRC of the data file is not equal to CRC subfield of catalog record" + $rq + "."
$found3 = $rng3.Find.Execute($target3)
if ($found3) {
  $xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' +
  '<w:p w14:paraId="4F9C12CE" w14:textId="77777777" w:rsidR="003B4B9F" w:rsidRDefault="003B4B9F">' +
  '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">1016: </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">' + $lq + 'CRC </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>values do not match</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>' + $rq + '.</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  $rng3.InsertXML($xml3)
} else {
  Write-Host "WARN: edit 3 target not found"
}

# ---------------------------------------------------------------------------
# Edit 4: after the "531: An error "bad cell name wrong agency"." paragraph,
# add a new paragraph: 1518b: "Producing Agency code is not a valid S-62 value"
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$target4 = "531: An error " + $lq + "bad cell name wrong agency" + $rq + "."
$found4 = $rng4.Find.Execute($target4)
if ($found4) {
  $xml4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' +
  '<w:p w14:paraId="16CAB07E" w14:textId="77777777" w:rsidR="003B4B9F" w:rsidRDefault="003B4B9F">' +
  '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">531: An error ' + $lq + 'bad cell name wrong </w:t></w:r>' +
  '<w:r><w:t>agency</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>' + $rq + '.</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>1518b: ' + $lq + 'Producing Agency code is not a valid S-62 value' + $rq + '</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  $rng4.InsertXML($xml4)
} else {
  Write-Host "WARN: edit 4 target not found"
}

Write-Host "Done."
